# Generate Report for Handback
#
# The localization-status workbook is regenerated for a handback: the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated for both language sheets (zh-cn, de-de) and the
# Overview/per-language Status cells flip from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook
$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column updates -------------------------------------------------
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# --- Hyperlink target file names (column I) + link relationships ----------
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a80759bc9d9809e88955ccef59d8e3f0273c454/e2e/d1e2288a-746e-40c3-9652-46756e81b961.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a80759bc9d9809e88955ccef59d8e3f0273c454/e2e/ed41a8cb-8751-48e6-8b81-9c53774b6ad5.md"
$mdName1 = "d1e2288a-746e-40c3-9652-46756e81b961.md"
$mdName2 = "ed41a8cb-8751-48e6-8b81-9c53774b6ad5.md"

function Set-HandbackRow($ws, $row, $mdUrl, $mdName, $xlfName, $handbackDate) {
    $targetCell = $ws.Range("I" + $row)
    $targetCell.Value = $mdName
    $ws.Hyperlinks.Add($targetCell, $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
    # match the look of the other hyperlink cells (A2/A3): underlined, cornflower blue
    $targetCell.Font.Underline = 2
    $targetCell.Font.Color = 15570276

    $ws.Range("J" + $row).Value = $xlfName
    $ws.Range("K" + $row).Value = $handbackDate
}

Set-HandbackRow $zh 2 $mdUrl1 $mdName1 "d1e2288a-746e-40c3-9652-46756e81b961.a25535389bff0d18ffcf3e1e51c228669ca8c574.zh-cn.xlf" "2016-08-25 06:28:36"
Set-HandbackRow $zh 3 $mdUrl2 $mdName2 "ed41a8cb-8751-48e6-8b81-9c53774b6ad5.0523956594df5eeaf05c8872e9f51acf069f1bf2.zh-cn.xlf" "2016-08-25 06:28:36"

Set-HandbackRow $de 2 $mdUrl1 $mdName1 "d1e2288a-746e-40c3-9652-46756e81b961.a25535389bff0d18ffcf3e1e51c228669ca8c574.de-de.xlf" "2016-08-25 06:28:43"
Set-HandbackRow $de 3 $mdUrl2 $mdName2 "ed41a8cb-8751-48e6-8b81-9c53774b6ad5.0523956594df5eeaf05c8872e9f51acf069f1bf2.de-de.xlf" "2016-08-25 06:28:43"

# --- Column width refresh (the longer text now needs more room) -----------
# COM ColumnWidth is stored with a fixed +5/6 character padding, so back it
# out to land on the same persisted <col width> the workbook shipped with.
$padding = 0.8333333333333334

$ov.Range("E1:F1").ColumnWidth = 29.9777047293527 - $padding

$zh.Range("C1").ColumnWidth = 29.9777047293527 - $padding
$zh.Range("I1:J1").ColumnWidth = 40 - $padding

$de.Range("C1").ColumnWidth = 29.9777047293527 - $padding
$de.Range("I1:J1").ColumnWidth = 40 - $padding

Write-Host "Handback report generated"
